$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = [double]"1"
$ws.Range("F2").Value = [double]"0.3333333333333333"
$ws.Range("G2").Value = [double]"0.1234603333333333"
$ws.Range("H2").Value = [double]"0.370381"
$ws.Range("I2").Value = [double]"0.002558470358543426"
$ws.Range("J2").Value = [double]"0.002636284444771545"
$ws.Range("M2").Value = [double]"103.7041626666667"
$ws.Range("N2").Value = [double]"311.112488"
$ws.Range("O2").Value = [double]"0.9879014414010097"
$ws.Range("P2").Value = [double]"0.9902837487998425"
$ws.Range("Q2").Value = [double]"12.80335049088089"
$ws.Range("R2").Value = [double]"115.230154417928"
$ws.Range("S2").Value = [double]"0.002527516554986809"
$ws.Range("T2").Value = [double]"0.002610669642871077"
$ws.Range("E3").Value = [double]"1"
$ws.Range("F3").Value = [double]"0.3333333333333333"
$ws.Range("G3").Value = [double]"0.1234603333333333"
$ws.Range("H3").Value = [double]"0.370381"
$ws.Range("I3").Value = [double]"0.002558470358543426"
$ws.Range("J3").Value = [double]"0.002636284444771545"
$ws.Range("O3").Value = [double]"0.003717513475672384"
$ws.Range("P3").Value = [double]"0.003726478195721577"
$ws.Range("Q3").Value = [double]"0.04817953086100001"
$ws.Range("R3").Value = [double]"0.4336157777490001"
$ws.Range("S3").Value = [double]"9.511148034993542E-06"
$ws.Range("T3").Value = [double]"9.824056501161129E-06"
$ws.Range("E4").Value = [double]"1"
$ws.Range("F4").Value = [double]"0.3333333333333333"
$ws.Range("G4").Value = [double]"0.1234603333333333"
$ws.Range("H4").Value = [double]"0.370381"
$ws.Range("I4").Value = [double]"0.002558470358543426"
$ws.Range("J4").Value = [double]"0.002636284444771545"
$ws.Range("M4").Value = [double]"0.09991866666666667"
$ws.Range("N4").Value = [double]"0.299756"
$ws.Range("O4").Value = [double]"0.0009518402375047094"
$ws.Range("P4").Value = [double]"0.0009541355839282337"
$ws.Range("Q4").Value = [double]"0.01233599189288889"
$ws.Range("R4").Value = [double]"0.111023927036"
$ws.Range("S4").Value = [double]"2.435255033724734E-06"
$ws.Range("T4").Value = [double]"2.515372798113018E-06"
$ws.Range("E5").Value = [double]"1"
$ws.Range("F5").Value = [double]"0.3333333333333333"
$ws.Range("G5").Value = [double]"0.1234603333333333"
$ws.Range("H5").Value = [double]"0.370381"
$ws.Range("I5").Value = [double]"0.002558470358543426"
$ws.Range("J5").Value = [double]"0.002636284444771545"
$ws.Range("M5").Value = [double]"0.7576035"
$ws.Range("N5").Value = [double]"1.515207"
$ws.Range("O5").Value = [double]"0.007217044816861706"
$ws.Range("P5").Value = [double]"0.004822965731185187"
$ws.Range("Q5").Value = [double]"0.0935339806445"
$ws.Range("R5").Value = [double]"0.561203883867"
$ws.Range("S5").Value = [double]"1.846459524022014E-05"
$ws.Range("T5").Value = [double]"1.271470953478973E-05"
$ws.Range("E6").Value = [double]"1"
$ws.Range("F6").Value = [double]"0.3333333333333333"
$ws.Range("G6").Value = [double]"0.1234603333333333"
$ws.Range("H6").Value = [double]"0.370381"
$ws.Range("I6").Value = [double]"0.002558470358543426"
$ws.Range("J6").Value = [double]"0.002636284444771545"
$ws.Range("K6").Value = [double]"1"
$ws.Range("L6").Value = [double]"0.3333333333333333"
$ws.Range("M6").Value = [double]"0.02227133333333333"
$ws.Range("N6").Value = [double]"0.066814"
$ws.Range("O6").Value = [double]"0.0002121600689515461"
$ws.Range("P6").Value = [double]"0.0002126716893225857"
$ws.Range("Q6").Value = [double]"0.002749626237111111"
$ws.Range("R6").Value = [double]"0.024746636134"
$ws.Range("S6").Value = [double]"5.428052476790602E-07"
$ws.Range("T6").Value = [double]"5.606630664044194E-07"
$ws.Range("I7").Value = [double]"0.9088918061291337"
$ws.Range("J7").Value = [double]"0.9365351146153917"
$ws.Range("M7").Value = [double]"103.7041626666667"
$ws.Range("N7").Value = [double]"311.112488"
$ws.Range("O7").Value = [double]"0.9879014414010097"
$ws.Range("P7").Value = [double]"0.9902837487998425"
$ws.Range("Q7").Value = [double]"4548.366297581846"
$ws.Range("R7").Value = [double]"40935.29667823661"
$ws.Range("S7").Value = [double]"0.8978955253525381"
$ws.Range("T7").Value = [double]"0.9274355041840203"
$ws.Range("I8").Value = [double]"0.9088918061291337"
$ws.Range("J8").Value = [double]"0.9365351146153917"
$ws.Range("O8").Value = [double]"0.003717513475672384"
$ws.Range("P8").Value = [double]"0.003726478195721577"
$ws.Range("S8").Value = [double]"0.003378817537213266"
$ws.Range("T8").Value = [double]"0.003489977684141866"
$ws.Range("I9").Value = [double]"0.9088918061291337"
$ws.Range("J9").Value = [double]"0.9365351146153917"
$ws.Range("M9").Value = [double]"0.09991866666666667"
$ws.Range("N9").Value = [double]"0.299756"
$ws.Range("O9").Value = [double]"0.0009518402375047094"
$ws.Range("P9").Value = [double]"0.0009541355839282337"
$ws.Range("Q9").Value = [double]"4.382338030410222"
$ws.Range("R9").Value = [double]"39.44104227369201"
$ws.Range("S9").Value = [double]"0.0008651197926120389"
$ws.Range("T9").Value = [double]"0.0008935814784528521"
$ws.Range("I10").Value = [double]"0.9088918061291337"
$ws.Range("J10").Value = [double]"0.9365351146153917"
$ws.Range("M10").Value = [double]"0.7576035"
$ws.Range("N10").Value = [double]"1.515207"
$ws.Range("O10").Value = [double]"0.007217044816861706"
$ws.Range("P10").Value = [double]"0.004822965731185187"
$ws.Range("Q10").Value = [double]"33.2277715544165"
$ws.Range("R10").Value = [double]"199.366629326499"
$ws.Range("S10").Value = [double]"0.006559512898512338"
$ws.Range("T10").Value = [double]"0.004516876763841626"
$ws.Range("I11").Value = [double]"0.9088918061291337"
$ws.Range("J11").Value = [double]"0.9365351146153917"
$ws.Range("K11").Value = [double]"1"
$ws.Range("L11").Value = [double]"0.3333333333333333"
$ws.Range("M11").Value = [double]"0.02227133333333333"
$ws.Range("N11").Value = [double]"0.066814"
$ws.Range("O11").Value = [double]"0.0002121600689515461"
$ws.Range("P11").Value = [double]"0.0002126716893225857"
$ws.Range("Q11").Value = [double]"0.9767995741997777"
$ws.Range("R11").Value = [double]"8.791196167798001"
$ws.Range("S11").Value = [double]"0.0001928305482578523"
$ws.Range("T11").Value = [double]"0.0001991745049351768"
$ws.Range("G12").Value = [double]"4.2730135"
$ws.Range("H12").Value = [double]"8.546027"
$ws.Range("I12").Value = [double]"0.08854972351232299"
$ws.Range("J12").Value = [double]"0.06082860093983664"
$ws.Range("M12").Value = [double]"103.7041626666667"
$ws.Range("N12").Value = [double]"311.112488"
$ws.Range("O12").Value = [double]"0.9879014414010097"
$ws.Range("P12").Value = [double]"0.9902837487998425"
$ws.Range("Q12").Value = [double]"443.1292870808627"
$ws.Range("R12").Value = [double]"2658.775722485176"
$ws.Range("S12").Value = [double]"0.08747839949348475"
$ws.Range("T12").Value = [double]"0.06023757497295105"
$ws.Range("G13").Value = [double]"4.2730135"
$ws.Range("H13").Value = [double]"8.546027"
$ws.Range("I13").Value = [double]"0.08854972351232299"
$ws.Range("J13").Value = [double]"0.06082860093983664"
$ws.Range("O13").Value = [double]"0.003717513475672384"
$ws.Range("P13").Value = [double]"0.003726478195721577"
$ws.Range("Q13").Value = [double]"1.6675136072805"
$ws.Range("R13").Value = [double]"10.005081643683"
$ws.Range("S13").Value = [double]"0.0003291847904241244"
$ws.Range("T13").Value = [double]"0.0002266764550785503"
$ws.Range("G14").Value = [double]"4.2730135"
$ws.Range("H14").Value = [double]"8.546027"
$ws.Range("I14").Value = [double]"0.08854972351232299"
$ws.Range("J14").Value = [double]"0.06082860093983664"
$ws.Range("M14").Value = [double]"0.09991866666666667"
$ws.Range("N14").Value = [double]"0.299756"
$ws.Range("O14").Value = [double]"0.0009518402375047094"
$ws.Range("P14").Value = [double]"0.0009541355839282337"
$ws.Range("Q14").Value = [double]"0.4269538115686667"
$ws.Range("R14").Value = [double]"2.561722869412"
$ws.Range("S14").Value = [double]"8.428518985894586E-05"
$ws.Range("T14").Value = [double]"5.803873267726854E-05"
$ws.Range("G15").Value = [double]"4.2730135"
$ws.Range("H15").Value = [double]"8.546027"
$ws.Range("I15").Value = [double]"0.08854972351232299"
$ws.Range("J15").Value = [double]"0.06082860093983664"
$ws.Range("M15").Value = [double]"0.7576035"
$ws.Range("N15").Value = [double]"1.515207"
$ws.Range("O15").Value = [double]"0.007217044816861706"
$ws.Range("P15").Value = [double]"0.004822965731185187"
$ws.Range("Q15").Value = [double]"3.23724998314725"
$ws.Range("R15").Value = [double]"12.948999932589"
$ws.Range("S15").Value = [double]"0.0006390673231091477"
$ws.Range("T15").Value = [double]"0.0002933742578087712"
$ws.Range("G16").Value = [double]"4.2730135"
$ws.Range("H16").Value = [double]"8.546027"
$ws.Range("I16").Value = [double]"0.08854972351232299"
$ws.Range("J16").Value = [double]"0.06082860093983664"
$ws.Range("K16").Value = [double]"1"
$ws.Range("L16").Value = [double]"0.3333333333333333"
$ws.Range("M16").Value = [double]"0.02227133333333333"
$ws.Range("N16").Value = [double]"0.066814"
$ws.Range("O16").Value = [double]"0.0002121600689515461"
$ws.Range("P16").Value = [double]"0.0002126716893225857"
$ws.Range("Q16").Value = [double]"0.09516570799633334"
$ws.Range("R16").Value = [double]"0.570994247978"
$ws.Range("S16").Value = [double]"1.878671544601479E-05"
$ws.Range("T16").Value = [double]"1.293652132100448E-05"
